{"js": "// Add a new \"Entorhinal\" row (50 | 62.5 | 66.66) to the end of the\n// accuracy table, right after the existing \"Frontal Cortex\" row.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\n\n// Office.js Table.addRows(location, count, values) \u2014 \"End\" appends after\n// the last existing row; values fills the new row's cells left to right.\ntable.addRows(\"End\", 1, [[\"Entorhinal\", \"50\", \"62.5\", \"66.66\"]]);\n\nawait context.sync();\n", "ps1": "# Add a new \"Entorhinal\" row (50 | 62.5 | 66.66) to the end of the\n# accuracy table, right after the existing \"Frontal Cortex\" row.\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n# Table.Rows.Add() with no \"before row\" argument appends a new row at the\n# end of the table, inheriting the formatting (e.g. trHeight) of the last row.\n$newRow = $table.Rows.Add()\n\n$newRow.Cells(1).Range.Text = \"Entorhinal\"\n$newRow.Cells(2).Range.Text = \"50\"\n$newRow.Cells(3).Range.Text = \"62.5\"\n$newRow.Cells(4).Range.Text = \"66.66\"\n"}
